$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.396.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.526.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.531.36"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0991"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("E11").Value = "  +1.48%  "

$ws.Range("E12").Value = "  -2.55%  "

$ws.Range("E13").Value = "  +1.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.973.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.309.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.505.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.421"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.68%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0767"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.05%  "

$ws.Range("E33").Value = "  +0.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.19%  "

$ws.Range("E35").Value = "  +0.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.20%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.54%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "284.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.601"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0924"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0509"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.61%  "

$ws.Range("E51").Value = "  -2.13%  "
